$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set cell values in the same order the original strings table was built ---
$ws.Range("L13").Value = "Feb. 21, 2017"
$ws.Range("J13").Value = "36 hr"
$ws.Range("J14").Value = "37 hr"
$ws.Range("G13").Value = "Done! [buggy]"
$ws.Range("I13").Value = "0.5650 / 0.6707 / 3.548"
$ws.Range("H13").Value = "0.6323 / 0.7717 / 1.902"
$ws.Range("H14").Value = "0.7340 / 0.8672 / 1.091"
$ws.Range("I14").Value = "0.5727 / 0.6958 / 3.692"
$ws.Range("K14").Value = "model is overfitting!"
$ws.Range("K13").Value = "After step 29,000 learning curves are `nbuggy, e.g., loss is NaN! Needs debugging, working on it. Reported results are the results after 29,000 updates right before the bug."
$ws.Range("A15").Value = "Repeating Experiment run_id 8 to confirm if the bug is repeatable."
$ws.Range("D15").Value = "EXP13.txt"
$ws.Range("B15").Value = "python -m basic.cli --mode train --noload --len_opt --cluster --TPRLSTMCell False --justTPR True --num_steps 40000 --num_epochs 24 --batch_size 40 --run_id 10 |& tee /home/hpalangi/QA/TPR_Stuff/Codes/TPR_ver1.0/Log_Files/EXP13.txt"

$ws.Range("L14").Value = "Feb. 21, 2017"
$ws.Range("G14").Value = "Done!"
$ws.Range("C15").Value = "DLT2 / 2"
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = 0
$ws.Range("K15").Value = "After step 29,000 learning curves are `nbuggy, e.g., loss is NaN! Needs debugging, working on it. Reported results are the results after 29,000 updates right before the bug."
$ws.Range("L15").Value = "Feb. 21, 2017"

# --- Formatting ---
$g13 = $ws.Range("G13")
$g13.Characters(7, 7).Font.Color = 255

$k13 = $ws.Range("K13")
$k13.WrapText = $true
$k13.Font.Color = 255
$k13.Interior.Color = 5296274

$ws.Rows.Item(15).RowHeight = 75

$a15 = $ws.Range("A15")
$a15.WrapText = $true
$a15.Interior.Color = 255

$ws.Range("B15").Interior.Color = 255
$ws.Range("C15").Interior.Color = 255
$ws.Range("D15").Interior.Color = 255
$ws.Range("E15").Interior.Color = 255
$ws.Range("F15").Interior.Color = 255

$k15 = $ws.Range("K15")
$k15.WrapText = $true
$k15.Font.Color = 255
$k15.Interior.Color = 255

$ws.Range("L15").Interior.Color = 255

# --- View state ---
$ws.Range("B15").Select()
